# Apply updated TPM-derived values to the NATMI LR-pairs sheet (Ntng2-Lrrc4)
# Values below come from re-running the scripts with the new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 6.840962666666667
    "H2" = 20.522888
    "I2" = 0.3877069165303996
    "J2" = 0.3877069165303996
    "M2" = 2.411650333333333
    "N2" = 7.234951
    "O2" = 0.2281979581963742
    "P2" = 0.2281979581963742
    "Q2" = 16.49800989538755
    "R2" = 148.482089058488
    "S2" = 0.08847392673084926
    "T2" = 0.08847392673084928
    "G3" = 6.840962666666667
    "H3" = 20.522888
    "I3" = 0.3877069165303996
    "J3" = 0.3877069165303996
    "O3" = 0.3161026464392645
    "P3" = 0.3161026464392646
    "Q3" = 22.85324825047467
    "R3" = 205.679234254272
    "S3" = 0.1225551823580664
    "T3" = 0.1225551823580664
    "G4" = 6.840962666666667
    "H4" = 20.522888
    "I4" = 0.3877069165303996
    "J4" = 0.3877069165303996
    "M4" = 2.855816333333333
    "N4" = 8.567449
    "O4" = 0.2702263455207323
    "P4" = 0.2702263455207323
    "Q4" = 19.53653291919022
    "R4" = 175.828796272712
    "S4" = 0.1047686231871215
    "T4" = 0.1047686231871215
    "G5" = 6.840962666666667
    "H5" = 20.522888
    "I5" = 0.3877069165303996
    "J5" = 0.3877069165303996
    "M5" = 1.960123333333333
    "N5" = 5.88037
    "O5" = 0.1854730498436289
    "P5" = 0.1854730498436289
    "Q5" = 13.40913054539556
    "R5" = 120.68217490856
    "S5" = 0.07190918425436249
    "T5" = 0.07190918425436249
    "I6" = 0.3116217912463337
    "J6" = 0.3116217912463337
    "M6" = 2.411650333333333
    "N6" = 7.234951
    "O6" = 0.2281979581963742
    "P6" = 0.2281979581963742
    "Q6" = 13.26037575395511
    "R6" = 119.343381785596
    "S6" = 0.0711114564919101
    "T6" = 0.07111145649191011
    "I7" = 0.3116217912463337
    "J7" = 0.3116217912463337
    "O7" = 0.3161026464392645
    "P7" = 0.3161026464392646
    "S7" = 0.09850447290111011
    "T7" = 0.09850447290111013
    "I8" = 0.3116217912463337
    "J8" = 0.3116217912463337
    "M8" = 2.855816333333333
    "N8" = 8.567449
    "O8" = 0.2702263455207323
    "P8" = 0.2702263455207323
    "Q8" = 15.70260710720044
    "R8" = 141.323463964804
    "S8" = 0.08420841783312129
    "T8" = 0.08420841783312129
    "I9" = 0.3116217912463337
    "J9" = 0.3116217912463337
    "M9" = 1.960123333333333
    "N9" = 5.88037
    "O9" = 0.1854730498436289
    "P9" = 0.1854730498436289
    "Q9" = 10.77767019739111
    "R9" = 96.99903177652
    "S9" = 0.05779744402019218
    "T9" = 0.05779744402019218
    "G10" = 2.722503999999999
    "H10" = 8.167511999999999
    "I10" = 0.1542960665791791
    "J10" = 0.1542960665791791
    "M10" = 2.411650333333333
    "N10" = 7.234951
    "O10" = 0.2281979581963742
    "P10" = 0.2281979581963742
    "Q10" = 6.565727679101331
    "R10" = 59.09154911191199
    "S10" = 0.03521004735110048
    "T10" = 0.03521004735110049
    "G11" = 2.722503999999999
    "H11" = 8.167511999999999
    "I11" = 0.1542960665791791
    "J11" = 0.1542960665791791
    "O11" = 0.3161026464392645
    "P11" = 0.3161026464392646
    "Q11" = 9.094927542591998
    "R11" = 81.85434788332799
    "S11" = 0.04877339498084748
    "T11" = 0.04877339498084748
    "G12" = 2.722503999999999
    "H12" = 8.167511999999999
    "I12" = 0.1542960665791791
    "J12" = 0.1542960665791791
    "M12" = 2.855816333333333
    "N12" = 8.567449
    "O12" = 0.2702263455207323
    "P12" = 0.2702263455207323
    "Q12" = 7.774971390765331
    "R12" = 69.97474251688799
    "S12" = 0.04169486219991518
    "T12" = 0.04169486219991518
    "G13" = 2.722503999999999
    "H13" = 8.167511999999999
    "I13" = 0.1542960665791791
    "J13" = 0.1542960665791791
    "M13" = 1.960123333333333
    "N13" = 5.88037
    "O13" = 0.1854730498436289
    "P13" = 0.1854730498436289
    "Q13" = 5.336443615493333
    "R13" = 48.02799253943999
    "S13" = 0.02861776204731598
    "T13" = 0.02861776204731598
    "G14" = 2.582743333333333
    "H14" = 7.74823
    "I14" = 0.1463752256440876
    "J14" = 0.1463752256440876
    "M14" = 2.411650333333333
    "N14" = 7.234951
    "O14" = 0.2281979581963742
    "P14" = 0.2281979581963742
    "Q14" = 6.228673820747777
    "R14" = 56.05806438673
    "S14" = 0.03340252762251434
    "T14" = 0.03340252762251435
    "G15" = 2.582743333333333
    "H15" = 7.74823
    "I15" = 0.1463752256440876
    "J15" = 0.1463752256440876
    "O15" = 0.3161026464392645
    "P15" = 0.3161026464392646
    "Q15" = 8.628036351013334
    "R15" = 77.65232715912001
    "S15" = 0.04626959619924059
    "T15" = 0.0462695961992406
    "G16" = 2.582743333333333
    "H16" = 7.74823
    "I16" = 0.1463752256440876
    "J16" = 0.1463752256440876
    "M16" = 2.855816333333333
    "N16" = 8.567449
    "O16" = 0.2702263455207323
    "P16" = 0.2702263455207323
    "Q16" = 7.375840596141111
    "R16" = 66.38256536527
    "S16" = 0.03955444230057438
    "T16" = 0.03955444230057438
    "G17" = 2.582743333333333
    "H17" = 7.74823
    "I17" = 0.1463752256440876
    "J17" = 0.1463752256440876
    "M17" = 1.960123333333333
    "N17" = 5.88037
    "O17" = 0.1854730498436289
    "P17" = 0.1854730498436289
    "Q17" = 5.062495471677778
    "R17" = 45.5624592451
    "S17" = 0.02714865952175829
    "T17" = 0.02714865952175829
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
